$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.596.43'
$ws.Range("E2").Value = '  +1.45%  '

$ws.Range("D3").Value = '''3.159.32'
$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''592.05'
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").Value = '''147.32'
$ws.Range("E6").Value = '  +0.94%  '

$ws.Range("D8").Value = '''3.156.12'
$ws.Range("E8").Value = '  +1.02%  '

$ws.Range("D9").Value = '''0.532'
$ws.Range("E9").Value = '  -0.76%  '

$ws.Range("E10").Value = '  -0.51%  '

$ws.Range("D11").Value = '''5.98'
$ws.Range("E11").Value = '  +4.95%  '

$ws.Range("D12").Value = '''0.464'
$ws.Range("E12").Value = '  -0.92%  '

$ws.Range("D13").Value = '''0.0000250'
$ws.Range("E13").Value = '  -2.08%  '

$ws.Range("D14").Value = '''37.32'
$ws.Range("E14").Value = '  +3.23%  '

$ws.Range("D15").Value = '''3.680.24'
$ws.Range("E15").Value = '  +0.89%  '

$ws.Range("E16").Value = '  -1.10%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '''64.337.53'
$ws.Range("E17").Value = '  +1.16%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '''7.24'
$ws.Range("E18").Value = '  +0.58%  '

$ws.Range("D19").Value = '''3.156.32'
$ws.Range("E19").Value = '  +0.96%  '

$ws.Range("D20").Value = '''471.02'
$ws.Range("E20").Value = '  +1.02%  '

$ws.Range("D21").Value = '''14.51'
$ws.Range("E21").Value = '  +1.72%  '

$ws.Range("D22").Value = '''0.738'
$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '''13.12'
$ws.Range("E24").Value = '  -1.20%  '

$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D25").Value = '''2.34'
$ws.Range("E25").Value = '  +8.10%  '

$ws.Range("D26").Value = '''81.60'
$ws.Range("E26").Value = '  -0.92%  '

$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("D28").Value = '''9.80'
$ws.Range("E28").Value = '  +12.13%  '

$ws.Range("E29").Value = '  +0.60%  '

$ws.Range("D30").Value = '''7.42'
$ws.Range("E30").Value = '  +8.60%  '

$ws.Range("D31").Value = '''2.24'
$ws.Range("E31").Value = '  +0.44%  '

$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("D33").Value = '''27.54'
$ws.Range("E33").Value = '  +1.59%  '

$ws.Range("E34").Value = '  +0.85%  '

$ws.Range("D35").Value = '''0.0₃0853'
$ws.Range("E35").Value = '  -1.01%  '

$ws.Range("D36").Value = '''1.07'
$ws.Range("E36").Value = '  +1.42%  '

$ws.Range("D37").Value = '''2.36'
$ws.Range("E37").Value = '  -1.33%  '

$ws.Range("D38").Value = '''6.14'
$ws.Range("E38").Value = '  +0.05%  '

$ws.Range("D39").Value = '''3.28'
$ws.Range("E39").Value = '  -2.23%  '

$ws.Range("D40").Value = '''51.93'
$ws.Range("E40").Value = '  +2.04%  '

$ws.Range("D41").Value = '''456.06'
$ws.Range("E41").Value = '  +1.74%  '

$ws.Range("D42").Value = '''9.17'
$ws.Range("E42").Value = '  +4.77%  '

$ws.Range("D43").Value = '''0.294'
$ws.Range("E43").Value = '  +5.92%  '

$ws.Range("D44").Value = '''0.0375'
$ws.Range("E44").Value = '  +0.74%  '

$ws.Range("D45").Value = '''2.942.16'
$ws.Range("E45").Value = '  +0.83%  '

$ws.Range("D46").Value = '''40.65'
$ws.Range("E46").Value = '  +15.88%  '

$ws.Range("E47").Value = '  -0.82%  '

$ws.Range("D48").Value = '''128.37'
$ws.Range("E48").Value = '  +2.58%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").Value = '''2.26'
$ws.Range("E50").Value = '  +3.13%  '

$ws.Range("E51").Value = '  -0.08%  '
